$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before BM (shifts old BM->BN, old BN->BO)
$ws.Columns("BM").Insert()

# New header cell for the inserted column: a timestamp, like the other snapshot columns
$ws.Range("BM1").Value = "2026-01-30 15:24:01"

# Rows 2-80: the new BM column mirrors the latest price already recorded in BL
$bmValues = @(45.92,169.95,169.95,169.95,179.95,179.95,179.95,339.95,619,659,659,749,749,809,809,809,809,809,849,899,899,909,909,909,969,969,969,969,969,999,999,1039,1039,1079,1079,1079,1079,1099,1099,1199,1219,1219,1219,1219,1219,1229,1229,1249,1249,1329,1329,1329,1329,1329,1329,1329,1349,1479,1479,1479,1579,1579,1579,1579,1579,1579,1579,1729,1729,1729,1829,1829,1829,1979,1979,1979,2479,2479,2479)
for ($i = 0; $i -lt $bmValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 65).Value = $bmValues[$i]
}

# Rows 81-206: no price snapshot yet for these products, so the new BM column is blank
for ($row = 81; $row -le 206; $row++) {
    $ws.Cells.Item($row, 65).Value = ""
}
